$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0.032
$ws.Range("E3").Value = 0.878
$ws.Range("F3").Value = 0.022
$ws.Range("G3").Value = 0.096
$ws.Range("H3").Value = 0.014
$ws.Range("I3").Value = 0.018
$ws.Range("D4").Value = 0.5
$ws.Range("E4").Value = 0.484
$ws.Range("F4").Value = 0.51
$ws.Range("G4").Value = 0.474
$ws.Range("H4").Value = 0.499
$ws.Range("I4").Value = 0.514
$ws.Range("D5").Value = 0.951
$ws.Range("E5").Value = 0.122
$ws.Range("F5").Value = 0.987
$ws.Range("G5").Value = 0.882
$ws.Range("H5").Value = 0.978
$ws.Range("I5").Value = 0.966
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.012
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.992
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.001
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("E10").Value = 0
$ws.Range("E14").Value = 0.447
$ws.Range("F14").Value = 0.597
$ws.Range("H14").Value = 0.348
$ws.Range("D15").Value = 0.511
$ws.Range("E15").Value = 0.494
$ws.Range("F15").Value = 0.495
$ws.Range("G15").Value = 0.513
$ws.Range("H15").Value = 0.507
$ws.Range("I15").Value = 0.501
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0.529
$ws.Range("F16").Value = 0.443
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0.656
$ws.Range("I16").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0.539
$ws.Range("F17").Value = 0.358
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0.788
$ws.Range("I17").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.568
$ws.Range("F18").Value = 0.316
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0.867
$ws.Range("I18").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0.59
$ws.Range("F19").Value = 0.277
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0.905
$ws.Range("E20").Value = 0.627
$ws.Range("F20").Value = 0.266
$ws.Range("H20").Value = 0.943
$ws.Range("E21").Value = 0.646
$ws.Range("F21").Value = 0.229
$ws.Range("H21").Value = 0.977
$ws.Range("E22").Value = 0.658
$ws.Range("F22").Value = 0.22
$ws.Range("H22").Value = 0.978
$ws.Range("D25").Value = 0.186
$ws.Range("E25").Value = 0.786
$ws.Range("F25").Value = 0.832
$ws.Range("G25").Value = 0.19
$ws.Range("H25").Value = 0.455
$ws.Range("I25").Value = 0.251
$ws.Range("D26").Value = 0.511
$ws.Range("E26").Value = 0.478
$ws.Range("F26").Value = 0.491
$ws.Range("G26").Value = 0.507
$ws.Range("H26").Value = 0.502
$ws.Range("I26").Value = 0.507
$ws.Range("D27").Value = 0.811
$ws.Range("E27").Value = 0.209
$ws.Range("F27").Value = 0.209
$ws.Range("G27").Value = 0.819
$ws.Range("H27").Value = 0.548
$ws.Range("I27").Value = 0.76
$ws.Range("D28").Value = 0.964
$ws.Range("E28").Value = 0.055
$ws.Range("F28").Value = 0.049
$ws.Range("G28").Value = 0.956
$ws.Range("H28").Value = 0.583
$ws.Range("I28").Value = 0.921
$ws.Range("D29").Value = 0.998
$ws.Range("E29").Value = 0.011
$ws.Range("F29").Value = 0.009
$ws.Range("G29").Value = 0.996
$ws.Range("H29").Value = 0.652
$ws.Range("I29").Value = 0.977
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0.999
$ws.Range("H30").Value = 0.681
$ws.Range("I30").Value = 0.997
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 0.739
$ws.Range("I31").Value = 0.998
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 0.784
$ws.Range("I32").Value = 1
$ws.Range("E33").Value = 0
$ws.Range("H33").Value = 0.825
$ws.Range("D36").Value = 0.116
$ws.Range("E36").Value = 0.025
$ws.Range("F36").Value = 0.842
$ws.Range("G36").Value = 0.071
$ws.Range("H36").Value = 0.996
$ws.Range("I36").Value = 0.099
$ws.Range("D37").Value = 0.491
$ws.Range("E37").Value = 0.479
$ws.Range("F37").Value = 0.49
$ws.Range("G37").Value = 0.499
$ws.Range("H37").Value = 0.49
$ws.Range("I37").Value = 0.487
$ws.Range("D38").Value = 0.907
$ws.Range("E38").Value = 0.966
$ws.Range("F38").Value = 0.174
$ws.Range("G38").Value = 0.933
$ws.Range("H38").Value = 0.004
$ws.Range("I38").Value = 0.907
$ws.Range("D39").Value = 0.992
$ws.Range("E39").Value = 1
$ws.Range("F39").Value = 0.016
$ws.Range("G39").Value = 0.999
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0.994
$ws.Range("D40").Value = 1
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = 0.001
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 1
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 1
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 1
$ws.Range("D42").Value = 1
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 1
$ws.Range("I42").Value = 1
